$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.300926208496094
$ws.Range("B1").Value = 1.678855299949646
$ws.Range("C1").Value = 2.285144090652466
$ws.Range("D1").Value = 6.256363868713379
$ws.Range("E1").Value = 2.802734851837158
